$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column values are stored as text (e.g. "270.06") even though they look
# numeric, so a leading apostrophe forces Excel to keep them as text instead
# of silently converting them to real numbers.
$ws.Range("D2").Value = "'269.63"
$ws.Range("D3").Value = "'22.83"
$ws.Range("D4").Value = "'6.337"
$ws.Range("D5").Value = "'0.06175"
$ws.Range("D6").Value = "'3.643"
$ws.Range("D7").Value = "'6.668"
$ws.Range("D8").Value = "'1.390"
$ws.Range("E8").Value = "7FTXTokenFTTWorstin24h"
$ws.Range("D9").Value = "'0.8294"
$ws.Range("D10").Value = "'0.01369"
$ws.Range("D12").Value = "'0.08315"
$ws.Range("D13").Value = "'0.03548"
$ws.Range("D14").Value = "'0.03204"
$ws.Range("D15").Value = "'0.09326"
$ws.Range("D16").Value = "'3.886"
$ws.Range("D17").Value = "'0.001643"
$ws.Range("D19").Value = "'0.006351"
$ws.Range("D20").Value = "'0.005671"
$ws.Range("D21").Value = "'0.001077"
$ws.Range("D26").Value = "'0.1239"
$ws.Range("D40").Value = "'0.04718"
$ws.Range("D41").Value = "'0.006948"
$ws.Range("D42").Value = "'0.003799"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("D43").Value = "'0.1161"
$ws.Range("D44").Value = "'0.01183"
$ws.Range("D45").Value = "'0.00006258"
$ws.Range("D48").Value = "'0.9199"
$ws.Range("D49").Value = "'0.002534"
$ws.Range("D50").Value = "'0.00002300"
